$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.028.49"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "1.678.46"
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'215.85"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("E6").Value = "  -2.67%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'0.255"
$ws.Range("E8").Value = "  +1.91%  "

$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'21.39"
$ws.Range("E9").Value = "  +5.51%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.0626"
$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("D11").Value = "'0.0888"
$ws.Range("E11").Value = "  -1.05%  "

$ws.Range("D12").Value = "1.915.14"
$ws.Range("E12").Value = "  +0.75%  "

$ws.Range("D13").Value = "1.663.40"
$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("E14").Value = "  +0.66%  "

$ws.Range("D15").Value = "'0.534"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("D16").Value = "'66.50"

$ws.Range("D17").Value = "27.016.08"
$ws.Range("E17").Value = "  +0.39%  "

$ws.Range("D18").Value = "'8.15"
$ws.Range("E18").Value = "  +1.77%  "

$ws.Range("D19").Value = "'235.67"
$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("D20").Value = "0.0₃0738"
$ws.Range("E20").Value = "  +0.86%  "

$ws.Range("E22").Value = "  +1.68%  "

$ws.Range("E23").Value = "  +1.28%  "

$ws.Range("E24").Value = "  -4.12%  "

$ws.Range("D25").Value = "'146.63"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("D26").Value = "'7.26"
$ws.Range("E26").Value = "  +1.84%  "

$ws.Range("D27").Value = "'16.42"
$ws.Range("E27").Value = "  +3.14%  "

$ws.Range("E28").Value = "  -2.21%  "

$ws.Range("E29").Value = "  +0.30%  "

$ws.Range("D30").Value = "'0.0498"
$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").Value = "  +0.36%  "

$ws.Range("D33").Value = "1.539.31"
$ws.Range("E33").Value = "  +5.54%  "

$ws.Range("E34").Value = "  +1.11%  "

$ws.Range("E35").Value = "  +4.63%  "

$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("E37").Value = "  +2.16%  "

$ws.Range("E38").Value = "  +1.84%  "

$ws.Range("E39").Value = "  +3.24%  "

$ws.Range("E40").Value = "  +6.30%  "

$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").Value = "'68.00"
$ws.Range("E42").Value = "  +3.20%  "

$ws.Range("D43").Value = "'5.60"
$ws.Range("E43").Value = "  -2.48%  "

$ws.Range("E44").Value = "  -0.60%  "

$ws.Range("D45").Value = "1.820.11"
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("D46").Value = "'0.781"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").Value = "'1.54"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("E49").Value = "  +2.41%  "

$ws.Range("D50").Value = "'8.01"
$ws.Range("E50").Value = "  +6.02%  "

$ws.Range("E51").Value = "  +0.00%  "
